$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before row 2 (shifts existing data down)
$ws.Range("A2:C8").Insert()
$ws.Range("A2:C8").ClearFormats()

# Fill newly inserted rows 2-8
$ws.Range("A2").Value = -0.1724167168140411
$ws.Range("B2").Value = -0.3089450895786285
$ws.Range("C2").Value = 0.9990701079368592
$ws.Range("A3").Value = -0.1458440721035003
$ws.Range("B3").Value = -0.0731511116027832
$ws.Range("C3").Value = -0.1902845203876495
$ws.Range("A4").Value = 0.0277943685650825
$ws.Range("B4").Value = -0.0499382354319095
$ws.Range("C4").Value = 0.04505131021142
$ws.Range("A5").Value = -0.042302418500185
$ws.Range("B5").Value = -0.052381694316864
$ws.Range("C5").Value = -0.0262672062963247
$ws.Range("A6").Value = 0.0780380368232727
$ws.Range("B6").Value = -0.0073303831741213
$ws.Range("C6").Value = 0.0215329993516206
$ws.Range("A7").Value = -0.0259617734700441
$ws.Range("B7").Value = -0.0493273697793483
$ws.Range("C7").Value = -0.0320704244077205
$ws.Range("A8").Value = -0.117286130785942
$ws.Range("B8").Value = -0.0560468845069408
$ws.Range("C8").Value = -0.0174096599221229

# Append 3 new rows at the end (29-31)
$ws.Range("A29").Value = 0.1461495161056518
$ws.Range("B29").Value = 0.1336267739534378
$ws.Range("C29").Value = -0.1892155110836029
$ws.Range("A30").Value = 0.117286130785942
$ws.Range("B30").Value = 0.7583891749382019
$ws.Range("C30").Value = 0.076052725315094
$ws.Range("A31").Value = 0.1061378344893455
$ws.Range("B31").Value = 0.2086104750633239
$ws.Range("C31").Value = -0.1314887404441833
